$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Training Data" — update traffic_volume (column D) values ---
$wsTrain = $wb.Worksheets.Item("Training Data")

$trainValues = @{
    2  = 21126
    3  = 16369
    4  = 17285
    5  = 16902
    6  = 14564
    7  = 15070
    8  = 19735
    9  = 16932
    10 = 16926
    11 = 16872
    12 = 16044
    13 = 21943
    14 = 19896
    15 = 14238
    16 = 17701
    17 = 16872
    18 = 16515
    19 = 17023
    20 = 20654
    21 = 17769
    22 = 16973
    23 = 17381
    24 = 17060
    25 = 23551
    26 = 21070
    27 = 17231
    28 = 14021
    29 = 11300
    30 = 13737
    31 = 14148
    32 = 15293
    33 = 16223
    34 = 16321
    35 = 18358
    36 = 17426
    37 = 22804
    38 = 20469
    39 = 14674
    40 = 13928
    41 = 15809
    42 = 17162
    43 = 31092
    44 = 19472
    45 = 18334
    46 = 17440
    47 = 18954
    48 = 13030
    49 = 23387
    50 = 21078
    51 = 16260
    52 = 19118
    53 = 17749
    54 = 17677
    55 = 22301
    56 = 18942
    57 = 18846
    58 = 20070
    59 = 18082
    60 = 25238
}

foreach ($row in $trainValues.Keys) {
    $wsTrain.Cells.Item($row, 4).Value = $trainValues[$row]
}

# --- Sheet 2: "Testing Data" — clear traffic_volume (column D) content ---
$wsTest = $wb.Worksheets.Item("Testing Data")

for ($row = 2; $row -le 13; $row++) {
    $wsTest.Cells.Item($row, 4).ClearContents()
}
